$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 20 - a clone of row 10's pattern (AutoPay All Fields / v2.7) with its own ID cell.
$ws.Range("A20").Value = "AutoPay All Fields"
$ws.Range("C20").Value = "19"
$ws.Range("D20").Value = "2.7"
$ws.Range("E20").Value = "10.50"
$ws.Range("G20").Value = "AutoPay"
$ws.Range("H20").Value = "en_US"
$ws.Range("I20").Value = "Jonty"
$ws.Range("J20").Value = "Smith"
$ws.Range("K20").Value = "15 Elm St"
$ws.Range("L20").Value = "Suite 600"
$ws.Range("M20").Value = "840"
$ws.Range("N20").Value = "Gambrills"
$ws.Range("O20").Value = "MD"
$ws.Range("P20").Value = "21054"
$ws.Range("S20").Value = "iahmed@govolution.com"
$ws.Range("S17").Copy()
$ws.Range("S20").PasteSpecial(-4122)
$ws.Range("T20").Value = "udf data 1"
$ws.Range("U20").Value = "udf data 2"
$ws.Range("V20").Value = "udf data 3"
$ws.Range("W20").Value = "udf data 4"
$ws.Range("X20").Value = "udf data 5"
$ws.Range("Y20").Value = "udf data 6"
$ws.Range("Z20").Value = "Orange"
$ws.Range("AA20").Value = "Soccer"
$ws.Range("AB20").Value = "udf data 9"
$ws.Range("AC20").Value = "udf data 10"
$ws.Range("AF20").Value = "|1234~Whole Wheat~$5| "

# Move the view the same way the author's session ended up.
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("AC23").Select()
